$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 460.07693
$ws.Range("I8").Value = 9.181818
$ws.Range("J8").Value = 2940
$ws.Range("K8").Value = 27.545454
$ws.Range("L8").Value = 8820
$ws.Range("M8").Value = 111.454546
$ws.Range("N8").Value = -9098
$ws.Range("H9").Value = 1162.375
$ws.Range("I9").Value = 400
$ws.Range("J9").Value = 1619.8
$ws.Range("K9").Value = 400
$ws.Range("L9").Value = 1619.8
$ws.Range("M9").Value = -231
$ws.Range("N9").Value = -1957.8
$ws.Range("H39").Value = 216.90909
$ws.Range("I39").Value = 179.11765
$ws.Range("K39").Value = 537.35295
$ws.Range("M39").Value = -241.35295
$ws.Range("H103").Value = 944.05884
$ws.Range("J103").Value = 1000.63635
$ws.Range("L103").Value = 3001.90905
$ws.Range("N103").Value = -4173.90905
$ws.Range("H110").Value = 33400.8
$ws.Range("J110").Value = 33400.8
$ws.Range("L110").Value = 33400.8
$ws.Range("N110").Value = -41580.8
$ws.Range("H138").Value = 3484.1072
$ws.Range("I138").Value = 4903.077
$ws.Range("J138").Value = 2254.3333
$ws.Range("K138").Value = 14709.231
$ws.Range("L138").Value = 6762.999899999999
$ws.Range("M138").Value = -9569.231
$ws.Range("N138").Value = -17042.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16017.17
$ws.Range("I32").Value = 6873.25
$ws.Range("J32").Value = 29950.762
$ws.Range("K32").Value = 6873.25
$ws.Range("L32").Value = 29950.762
$ws.Range("M32").Value = -6586.25
$ws.Range("N32").Value = -30524.762
$ws.Range("H45").Value = 1668.5
$ws.Range("I45").Value = 1337
$ws.Range("K45").Value = 1337
$ws.Range("M45").Value = -960
$ws.Range("H61").Value = 15535.842
$ws.Range("I61").Value = 4560.2
$ws.Range("J61").Value = 56694.5
$ws.Range("K61").Value = 4560.2
$ws.Range("L61").Value = 56694.5
$ws.Range("M61").Value = -4348.2
$ws.Range("N61").Value = -57118.5
$ws.Range("H122").Value = 3993.4
$ws.Range("I122").Value = 3707.5715
$ws.Range("K122").Value = 11122.7145
$ws.Range("M122").Value = -8672.7145
$ws.Range("H136").Value = 15535.842
$ws.Range("I136").Value = 4560.2
$ws.Range("J136").Value = 56694.5
$ws.Range("K136").Value = 13680.6
$ws.Range("L136").Value = 170083.5
$ws.Range("M136").Value = -11130.6
$ws.Range("N136").Value = -175183.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 15925.796
$ws.Range("I20").Value = 8158.92
$ws.Range("K20").Value = 8158.92
$ws.Range("M20").Value = -7911.92
$ws.Range("H33").Value = 8988.6
$ws.Range("I33").Value = 1973
$ws.Range("K33").Value = 1973
$ws.Range("M33").Value = -1637
$ws.Range("H38").Value = 19990
$ws.Range("J38").Value = 19990
$ws.Range("L38").Value = 19990
$ws.Range("N38").Value = -20822

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4030.2632
$ws.Range("I16").Value = 1179.909
$ws.Range("J16").Value = 7949.5
$ws.Range("K16").Value = 1179.909
$ws.Range("L16").Value = 7949.5
$ws.Range("M16").Value = -892.9090000000001
$ws.Range("N16").Value = -8523.5
$ws.Range("H22").Value = 1384.3846
$ws.Range("I22").Value = 599.4
$ws.Range("J22").Value = 1875
$ws.Range("K22").Value = 599.4
$ws.Range("L22").Value = 1875
$ws.Range("M22").Value = -249.4
$ws.Range("N22").Value = -2575
$ws.Range("H107").Value = 7029
$ws.Range("I107").Value = 3233
$ws.Range("J107").Value = 9137.888999999999
$ws.Range("K107").Value = 3233
$ws.Range("L107").Value = 9137.888999999999
$ws.Range("M107").Value = -1313
$ws.Range("N107").Value = -12977.889
$ws.Range("H113").Value = 4030.2632
$ws.Range("I113").Value = 1179.909
$ws.Range("J113").Value = 7949.5
$ws.Range("K113").Value = 1179.909
$ws.Range("L113").Value = 7949.5
$ws.Range("M113").Value = 990.0909999999999
$ws.Range("N113").Value = -12289.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2028.75
$ws.Range("I132").Value = 1501.4
$ws.Range("K132").Value = 13512.6
$ws.Range("M132").Value = -10982.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 6290.7954
$ws.Range("J15").Value = 6290.7954
$ws.Range("L15").Value = 6290.7954
$ws.Range("N15").Value = -6866.7954
$ws.Range("H70").Value = 12227.24
$ws.Range("I70").Value = 13545
$ws.Range("K70").Value = 13545
$ws.Range("M70").Value = -13275
$ws.Range("H73").Value = 12227.24
$ws.Range("I73").Value = 13545
$ws.Range("K73").Value = 13545
$ws.Range("M73").Value = -12609
$ws.Range("H81").Value = 6290.7954
$ws.Range("J81").Value = 6290.7954
$ws.Range("L81").Value = 6290.7954
$ws.Range("N81").Value = -8286.795399999999
$ws.Range("H84").Value = 6290.7954
$ws.Range("J84").Value = 6290.7954
$ws.Range("L84").Value = 18872.3862
$ws.Range("N84").Value = -28856.3862
$ws.Range("H122").Value = 2430.2727
$ws.Range("I122").Value = 1821
$ws.Range("J122").Value = 3496.5
$ws.Range("K122").Value = 5463
$ws.Range("L122").Value = 10489.5
$ws.Range("M122").Value = -3013
$ws.Range("N122").Value = -15389.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2594.7273
$ws.Range("I16").Value = 2614.3
$ws.Range("J16").Value = 2399
$ws.Range("K16").Value = 2614.3
$ws.Range("L16").Value = 2399
$ws.Range("M16").Value = -2444.3
$ws.Range("N16").Value = -2739
$ws.Range("H55").Value = 1982.9535
$ws.Range("I55").Value = 1004.2273
$ws.Range("J55").Value = 3008.2856
$ws.Range("K55").Value = 1004.2273
$ws.Range("L55").Value = 3008.2856
$ws.Range("M55").Value = -831.2273
$ws.Range("N55").Value = -3354.2856
$ws.Range("H61").Value = 2485.3333
$ws.Range("I61").Value = 1593.6666
$ws.Range("K61").Value = 1593.6666
$ws.Range("M61").Value = -1391.6666
$ws.Range("H113").Value = 2485.3333
$ws.Range("I113").Value = 1593.6666
$ws.Range("K113").Value = 1593.6666
$ws.Range("M113").Value = 576.3334
$ws.Range("H132").Value = 4458229.5
$ws.Range("I132").Value = 3437.4285
$ws.Range("K132").Value = 10312.2855
$ws.Range("M132").Value = -7782.2855
$ws.Range("H135").Value = 84983.336
$ws.Range("J135").Value = 84983.336
$ws.Range("L135").Value = 84983.336
$ws.Range("N135").Value = -95123.336
$ws.Range("H136").Value = 11281.415
$ws.Range("I136").Value = 9893.814
$ws.Range("J136").Value = 13957.5
$ws.Range("K136").Value = 29681.442
$ws.Range("L136").Value = 41872.5
$ws.Range("M136").Value = -27131.442
$ws.Range("N136").Value = -46972.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4501.3076
$ws.Range("I122").Value = 2208.9524
$ws.Range("J122").Value = 7175.722
$ws.Range("K122").Value = 6626.8572
$ws.Range("L122").Value = 21527.166
$ws.Range("M122").Value = -4176.8572
$ws.Range("N122").Value = -26427.166
$ws.Range("H132").Value = 20736.25
$ws.Range("I132").Value = 1577.2
$ws.Range("K132").Value = 4731.6
$ws.Range("M132").Value = -2201.6
